# Auto-generated by gen_ps1.py from the OOXML diff.
# Applies updated market-price / profit figures (static values, no formulas)
# to the leve-profit tracking sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 3756.84  # H64: 3778.682 -> 3756.84
$ws.Cells.Item(64, 10).Value = 3704.5  # J64: 3769.2 -> 3704.5
$ws.Cells.Item(64, 12).Value = 3704.5  # L64: 3769.2 -> 3704.5
$ws.Cells.Item(64, 14).Value = -4200.5  # N64: -4265.2 -> -4200.5
$ws.Cells.Item(67, 8).Value = 3756.84  # H67: 3778.682 -> 3756.84
$ws.Cells.Item(67, 10).Value = 3704.5  # J67: 3769.2 -> 3704.5
$ws.Cells.Item(67, 12).Value = 3704.5  # L67: 3769.2 -> 3704.5
$ws.Cells.Item(67, 14).Value = -5420.5  # N67: -5485.2 -> -5420.5
$ws.Cells.Item(74, 8).Value = 4728.933  # H74: 4842.154 -> 4728.933
$ws.Cells.Item(74, 9).Value = 4161.25  # I74: 4194.9 -> 4161.25
$ws.Cells.Item(74, 11).Value = 4161.25  # K74: 4194.9 -> 4161.25
$ws.Cells.Item(74, 13).Value = -3225.25  # M74: -3258.9 -> -3225.25
$ws.Cells.Item(77, 8).Value = 4728.933  # H77: 4842.154 -> 4728.933
$ws.Cells.Item(77, 9).Value = 4161.25  # I77: 4194.9 -> 4161.25
$ws.Cells.Item(77, 11).Value = 20806.25  # K77: 20974.5 -> 20806.25
$ws.Cells.Item(77, 13).Value = -16126.25  # M77: -16294.5 -> -16126.25
$ws.Cells.Item(80, 8).Value = 15626602  # H80: 13890379 -> 15626602
$ws.Cells.Item(80, 10).Value = 589.5  # J80: 592.8333 -> 589.5
$ws.Cells.Item(80, 12).Value = 1768.5  # L80: 1778.4999 -> 1768.5
$ws.Cells.Item(80, 14).Value = -3764.5  # N80: -3774.4999 -> -3764.5
$ws.Cells.Item(83, 8).Value = 15626602  # H83: 13890379 -> 15626602
$ws.Cells.Item(83, 10).Value = 589.5  # J83: 592.8333 -> 589.5
$ws.Cells.Item(83, 12).Value = 5305.5  # L83: 5335.4997 -> 5305.5
$ws.Cells.Item(83, 14).Value = -15289.5  # N83: -15319.4997 -> -15289.5
$ws.Cells.Item(92, 8).Value = 58823830  # H92: 62500228 -> 58823830
$ws.Cells.Item(92, 10).Value = 1749  # J92: 1998 -> 1749
$ws.Cells.Item(92, 12).Value = 1749  # L92: 1998 -> 1749
$ws.Cells.Item(92, 14).Value = -4245  # N92: -4494 -> -4245
$ws.Cells.Item(100, 8).Value = 1668.7273  # H100: 1402.9375 -> 1668.7273
$ws.Cells.Item(100, 9).Value = 1050.2858  # I100: 1158.8334 -> 1050.2858
$ws.Cells.Item(100, 10).Value = 2751  # J100: 1549.4 -> 2751
$ws.Cells.Item(100, 11).Value = 1050.2858  # K100: 1158.8334 -> 1050.2858
$ws.Cells.Item(100, 12).Value = 2751  # L100: 1549.4 -> 2751
$ws.Cells.Item(100, 13).Value = -509.2858000000001  # M100: -617.8334 -> -509.2858000000001
$ws.Cells.Item(100, 14).Value = -3833  # N100: -2631.4 -> -3833
$ws.Cells.Item(106, 8).Value = 3112.3572  # H106: 3259.4614 -> 3112.3572
$ws.Cells.Item(106, 9).Value = 3131.0833  # I106: 3306.6365 -> 3131.0833
$ws.Cells.Item(106, 11).Value = 3131.0833  # K106: 3306.6365 -> 3131.0833
$ws.Cells.Item(106, 13).Value = -2500.0833  # M106: -2675.6365 -> -2500.0833
$ws.Cells.Item(132, 8).Value = 7074.722  # H132: 7485.4707 -> 7074.722
$ws.Cells.Item(132, 9).Value = 1925.8064  # I132: 2052.276 -> 1925.8064
$ws.Cells.Item(132, 11).Value = 5777.4192  # K132: 6156.828 -> 5777.4192
$ws.Cells.Item(132, 13).Value = -3247.4192  # M132: -3626.828 -> -3247.4192
$ws.Cells.Item(138, 8).Value = 5681.618  # H138: 5710.655 -> 5681.618
$ws.Cells.Item(138, 9).Value = 2013.7142  # I138: 2116.3333 -> 2013.7142
$ws.Cells.Item(138, 10).Value = 6216.521  # J138: 6150.7754 -> 6216.521
$ws.Cells.Item(138, 11).Value = 6041.142599999999  # K138: 6348.999899999999 -> 6041.142599999999
$ws.Cells.Item(138, 12).Value = 18649.563  # L138: 18452.3262 -> 18649.563
$ws.Cells.Item(138, 13).Value = -901.1425999999992  # M138: -1208.999899999999 -> -901.1425999999992
$ws.Cells.Item(138, 14).Value = -28929.563  # N138: -28732.3262 -> -28929.563

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 11385.909  # H32: 10986.521 -> 11385.909
$ws.Cells.Item(32, 9).Value = 7166.2  # I32: 6855.8125 -> 7166.2
$ws.Cells.Item(32, 11).Value = 7166.2  # K32: 6855.8125 -> 7166.2
$ws.Cells.Item(32, 13).Value = -6879.2  # M32: -6568.8125 -> -6879.2
$ws.Cells.Item(74, 8).Value = 1405.4  # H74: 1021.95 -> 1405.4
$ws.Cells.Item(74, 9).Value = 1255  # I74: 912.5 -> 1255
$ws.Cells.Item(74, 11).Value = 1255  # K74: 912.5 -> 1255
$ws.Cells.Item(74, 13).Value = -381  # M74: -38.5 -> -381
$ws.Cells.Item(77, 8).Value = 1405.4  # H77: 1021.95 -> 1405.4
$ws.Cells.Item(77, 9).Value = 1255  # I77: 912.5 -> 1255
$ws.Cells.Item(77, 11).Value = 6275  # K77: 4562.5 -> 6275
$ws.Cells.Item(77, 13).Value = -1907  # M77: -194.5 -> -1907
$ws.Cells.Item(102, 8).Value = 1719.6364  # H102: 1750.8 -> 1719.6364
$ws.Cells.Item(102, 9).Value = 1611.7778  # I102: 1637.25 -> 1611.7778
$ws.Cells.Item(102, 11).Value = 1611.7778  # K102: 1637.25 -> 1611.7778
$ws.Cells.Item(102, 13).Value = 10.22219999999993  # M102: -15.25 -> 10.22219999999993
$ws.Cells.Item(132, 8).Value = 43253.188  # H132: 22065.469 -> 43253.188
$ws.Cells.Item(132, 9).Value = 64450.332  # I132: 23763.88 -> 64450.332
$ws.Cells.Item(132, 11).Value = 193350.996  # K132: 71291.64 -> 193350.996
$ws.Cells.Item(132, 13).Value = -190820.996  # M132: -68761.64 -> -190820.996

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1917.6  # H94: 2158.8 -> 1917.6
$ws.Cells.Item(94, 9).Value = 1796.5  # I94: 1799 -> 1796.5
$ws.Cells.Item(94, 10).Value = 1998.3334  # J94: 2248.75 -> 1998.3334
$ws.Cells.Item(94, 11).Value = 1796.5  # K94: 1799 -> 1796.5
$ws.Cells.Item(94, 12).Value = 1998.3334  # L94: 2248.75 -> 1998.3334
$ws.Cells.Item(94, 13).Value = -1345.5  # M94: -1348 -> -1345.5
$ws.Cells.Item(94, 14).Value = -2900.3334  # N94: -3150.75 -> -2900.3334
$ws.Cells.Item(99, 8).Value = 2136.5  # H99: 1243.8334 -> 2136.5
$ws.Cells.Item(99, 9).Value = 1699.5  # I99: 1083.1818 -> 1699.5
$ws.Cells.Item(99, 10).Value = 3010.5  # J99: 3011 -> 3010.5
$ws.Cells.Item(99, 11).Value = 1699.5  # K99: 1083.1818 -> 1699.5
$ws.Cells.Item(99, 12).Value = 3010.5  # L99: 3011 -> 3010.5
$ws.Cells.Item(99, 13).Value = -201.5  # M99: 414.8181999999999 -> -201.5
$ws.Cells.Item(99, 14).Value = -6006.5  # N99: -6007 -> -6006.5
$ws.Cells.Item(105, 8).Value = 2917.64  # H105: 2989.2083 -> 2917.64
$ws.Cells.Item(105, 9).Value = 2863.9048  # I105: 2947.1 -> 2863.9048
$ws.Cells.Item(105, 11).Value = 2863.9048  # K105: 2947.1 -> 2863.9048
$ws.Cells.Item(105, 13).Value = -1116.9048  # M105: -1200.1 -> -1116.9048
$ws.Cells.Item(134, 8).Value = 2250.8948  # H134: 2356.5278 -> 2250.8948
$ws.Cells.Item(134, 9).Value = 1743.1212  # I134: 1833.0322 -> 1743.1212
$ws.Cells.Item(134, 11).Value = 5229.363600000001  # K134: 5499.096600000001 -> 5229.363600000001
$ws.Cells.Item(134, 13).Value = -2694.363600000001  # M134: -2964.096600000001 -> -2694.363600000001

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(64, 8).Value = 39270.168  # H64: 39269.57 -> 39270.168
$ws.Cells.Item(64, 10).Value = 39270.168  # J64: 39269.57 -> 39270.168
$ws.Cells.Item(64, 12).Value = 39270.168  # L64: 39269.57 -> 39270.168
$ws.Cells.Item(64, 14).Value = -39766.168  # N64: -39765.57 -> -39766.168
$ws.Cells.Item(67, 8).Value = 39270.168  # H67: 39269.57 -> 39270.168
$ws.Cells.Item(67, 10).Value = 39270.168  # J67: 39269.57 -> 39270.168
$ws.Cells.Item(67, 12).Value = 39270.168  # L67: 39269.57 -> 39270.168
$ws.Cells.Item(67, 14).Value = -40986.168  # N67: -40985.57 -> -40986.168
$ws.Cells.Item(125, 8).Value = 15325  # H125: 0 -> 15325
$ws.Cells.Item(125, 10).Value = 15325  # J125: 0 -> 15325
$ws.Cells.Item(125, 12).Value = 15325  # L125: 0 -> 15325
$ws.Cells.Item(125, 14).Value = -20245  # N125: None -> -20245

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 2052841.9  # H4: 2113209.5 -> 2052841.9
$ws.Cells.Item(4, 9).Value = 61705.793  # I4: 63897.215 -> 61705.793
$ws.Cells.Item(4, 11).Value = 185117.379  # K4: 191691.645 -> 185117.379
$ws.Cells.Item(4, 13).Value = -185005.379  # M4: -191579.645 -> -185005.379
$ws.Cells.Item(37, 8).Value = 1000000000  # H37: 66775932 -> 1000000000
$ws.Cells.Item(37, 10).Value = 1000000000  # J37: 66775932 -> 1000000000
$ws.Cells.Item(37, 12).Value = 3000000000  # L37: 200327796 -> 3000000000
$ws.Cells.Item(37, 14).Value = -3000000224  # N37: -200328020 -> -3000000224
$ws.Cells.Item(38, 8).Value = 216  # H38: 253.6 -> 216
$ws.Cells.Item(38, 9).Value = 141.33333  # I38: 198 -> 141.33333
$ws.Cells.Item(38, 11).Value = 423.99999  # K38: 594 -> 423.99999
$ws.Cells.Item(38, 13).Value = -76.99998999999997  # M38: -247 -> -76.99998999999997
$ws.Cells.Item(52, 8).Value = 5488.4443  # H52: 5488.778 -> 5488.4443
$ws.Cells.Item(52, 10).Value = 5488.4443  # J52: 5488.778 -> 5488.4443
$ws.Cells.Item(52, 12).Value = 16465.3329  # L52: 16466.334 -> 16465.3329
$ws.Cells.Item(52, 14).Value = -16997.3329  # N52: -16998.334 -> -16997.3329
$ws.Cells.Item(107, 8).Value = 813  # H107: 954.6667 -> 813
$ws.Cells.Item(107, 10).Value = 813  # J107: 954.6667 -> 813
$ws.Cells.Item(107, 12).Value = 2439  # L107: 2864.0001 -> 2439
$ws.Cells.Item(107, 14).Value = -6279  # N107: -6704.0001 -> -6279
$ws.Cells.Item(114, 8).Value = 1744.5  # H114: 1845.4706 -> 1744.5
$ws.Cells.Item(114, 9).Value = 385.33334  # I114: 456.8 -> 385.33334
$ws.Cells.Item(114, 11).Value = 1156.00002  # K114: 1370.4 -> 1156.00002
$ws.Cells.Item(114, 13).Value = 2097.99998  # M114: 1883.6 -> 2097.99998
$ws.Cells.Item(118, 8).Value = 7054.8887  # H118: 6874.5 -> 7054.8887
$ws.Cells.Item(118, 10).Value = 8999.333000000001  # J118: 9250 -> 8999.333000000001
$ws.Cells.Item(118, 12).Value = 26997.999  # L118: 27750 -> 26997.999
$ws.Cells.Item(118, 14).Value = -29483.999  # N118: -30236 -> -29483.999

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(92, 8).Value = 7861.75  # H92: 21777.223 -> 7861.75
$ws.Cells.Item(92, 10).Value = 7861.75  # J92: 21777.223 -> 7861.75
$ws.Cells.Item(92, 12).Value = 7861.75  # L92: 21777.223 -> 7861.75
$ws.Cells.Item(92, 14).Value = -11605.75  # N92: -25521.223 -> -11605.75
$ws.Cells.Item(93, 14).Value = -35788.1  # N93: -38994 -> -35788.1
$ws.Cells.Item(93, 8).Value = 32044.1  # H93: 35250 -> 32044.1
$ws.Cells.Item(93, 10).Value = 32044.1  # J93: 35250 -> 32044.1
$ws.Cells.Item(93, 12).Value = 32044.1  # L93: 35250 -> 32044.1
$ws.Cells.Item(123, 8).Value = 36853  # H123: 38069 -> 36853
$ws.Cells.Item(123, 10).Value = 36853  # J123: 38069 -> 36853
$ws.Cells.Item(123, 12).Value = 36853  # L123: 38069 -> 36853
$ws.Cells.Item(123, 14).Value = -41753  # N123: -42969 -> -41753
$ws.Cells.Item(132, 8).Value = 4905.3125  # H132: 5041.857 -> 4905.3125
$ws.Cells.Item(132, 9).Value = 4362.154  # I132: 4437.1816 -> 4362.154
$ws.Cells.Item(132, 11).Value = 13086.462  # K132: 13311.5448 -> 13086.462
$ws.Cells.Item(132, 13).Value = -10556.462  # M132: -10781.5448 -> -10556.462

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 3683.1  # H46: 3836.2 -> 3683.1
$ws.Cells.Item(46, 10).Value = 3592.4075  # J46: 3670.4167 -> 3592.4075
$ws.Cells.Item(46, 12).Value = 3592.4075  # L46: 3670.4167 -> 3592.4075
$ws.Cells.Item(46, 14).Value = -3968.4075  # N46: -4046.4167 -> -3968.4075
$ws.Cells.Item(141, 8).Value = 85332.664  # H141: 88632.664 -> 85332.664
$ws.Cells.Item(141, 10).Value = 79000  # J141: 83950 -> 79000
$ws.Cells.Item(141, 12).Value = 79000  # L141: 83950 -> 79000
$ws.Cells.Item(141, 14).Value = -89360  # N141: -94310 -> -89360

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 1110.2727  # H100: 1133.6 -> 1110.2727
$ws.Cells.Item(100, 9).Value = 1083.5  # I100: 1083.75 -> 1083.5
$ws.Cells.Item(100, 10).Value = 1181.6666  # J100: 1333 -> 1181.6666
$ws.Cells.Item(100, 11).Value = 2167  # K100: 2167.5 -> 2167
$ws.Cells.Item(100, 12).Value = 2363.3332  # L100: 2666 -> 2363.3332
$ws.Cells.Item(100, 13).Value = -1626  # M100: -1626.5 -> -1626
$ws.Cells.Item(100, 14).Value = -3445.3332  # N100: -3748 -> -3445.3332
$ws.Cells.Item(108, 8).Value = 100626  # H108: 95313 -> 100626
$ws.Cells.Item(108, 10).Value = 100626  # J108: 95313 -> 100626
$ws.Cells.Item(108, 12).Value = 100626  # L108: 95313 -> 100626
$ws.Cells.Item(108, 14).Value = -108306  # N108: -102993 -> -108306
$ws.Cells.Item(132, 8).Value = 0  # H132: 499.5 -> 0
$ws.Cells.Item(132, 9).Value = 0  # I132: 499 -> 0
$ws.Cells.Item(132, 10).Value = 0  # J132: 500 -> 0
$ws.Cells.Item(132, 11).Value = 0  # K132: 1497 -> 0
$ws.Cells.Item(132, 12).Value = 0  # L132: 1500 -> 0
$ws.Cells.Item(132, 13).ClearContents()  # M132: was 1033, now blank
$ws.Cells.Item(132, 14).ClearContents()  # N132: was -6560, now blank

